# Atomix sheet: add Eyal's RX/TX perf numbers below the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atomix")
$ws.Activate()

# RX sub-table header (entered before the section title, matching shared-string order)
$ws.Cells.Item(66, 1).Value = "RX"
$ws.Cells.Item(66, 2).Value = "1-thread"
$ws.Cells.Item(66, 3).Value = "2-thread-optimal"
$ws.Cells.Item(66, 4).Value = "2-thread-old"

# RX data rows
$rxData = @(
    @(6,  71.962000000000003, 120.77,             96.225999999999999),
    @(9,  65.054000000000002, 110.726,            90.411000000000001),
    @(12, 54.067,             89.581999999999994, 83.775999999999996),
    @(18, 48.688000000000002, 74.278000000000006, 73.936999999999998),
    @(24, 36.741,             63.968000000000004, 62.466999999999999),
    @(36, 31.670999999999999, 50.856999999999999, 50.816000000000003),
    @(48, 28.651,             46.414999999999999, 42.634999999999998),
    @(54, 25.201000000000001, 50.433,             38.588000000000001)
)
$r = 67
foreach ($row in $rxData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
    $r++
}

# TX sub-table header (row 75 is intentionally left blank)
$ws.Cells.Item(76, 1).Value = "TX"
$ws.Cells.Item(76, 2).Value = "1-thread"
$ws.Cells.Item(76, 3).Value = "2-thread-optimal"
$ws.Cells.Item(76, 4).Value = "2-thread-old"

# TX data rows
$txData = @(
    @(6,  29.245000000000001, 35.978999999999999, 32.814),
    @(9,  43.296999999999997, 56.683,              50.832000000000001),
    @(12, 51.829000000000001, 74.492999999999995,  60.854999999999997),
    @(18, 74.122,             110.402,             96.483000000000004),
    @(24, 94.540999999999997, 149.12799999999999,  125.709),
    @(36, 108.226,            197.30799999999999,  189.06700000000001),
    @(48, 109.64100000000001, 193.47200000000001,  191.03299999999999),
    @(54, 125.764,            205.67599999999999,  195.65299999999999)
)
$r = 77
foreach ($row in $txData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 8).Value = $row[3]
    $r++
}

# Section title (bold, like the other row labels on this sheet) - added last
$ws.Cells.Item(65, 1).Value = "Eyal's numbers:"
$ws.Cells.Item(65, 1).Font.Bold = $true

# Match the author's final on-screen selection over the new TX "2-thread-optimal" column
$ws.Range("D77:D84").Select()
